# Auto-generated edit script applying the Famfrit_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 73.333336
$ws.Range("I5").Value = 73.333336
$ws.Range("K5").Value = 73.333336
$ws.Range("M5").Value = 41.666664
$ws.Range("H12").Value = 251.66667
$ws.Range("I12").Value = 336.66666
$ws.Range("J12").Value = 166.66667
$ws.Range("K12").Value = 336.66666
$ws.Range("L12").Value = 166.66667
$ws.Range("M12").Value = -166.66666
$ws.Range("N12").Value = -506.66667
$ws.Range("H96").Value = 37988.96
$ws.Range("I96").Value = 51651.39
$ws.Range("J96").Value = 2857
$ws.Range("K96").Value = 154954.17
$ws.Range("L96").Value = 8571
$ws.Range("M96").Value = -153581.17
$ws.Range("N96").Value = -11317
$ws.Range("H135").Value = 3534.2942
$ws.Range("I135").Value = 2033
$ws.Range("J135").Value = 4585.2
$ws.Range("K135").Value = 18297
$ws.Range("L135").Value = 41266.8
$ws.Range("M135").Value = -15762
$ws.Range("N135").Value = -46336.8
$ws.Range("H138").Value = 9213.68
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 9213.68
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 27641.04
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -37921.04
$ws.Range("H141").Value = 14394.25
$ws.Range("I141").Value = 30208.6
$ws.Range("K141").Value = 90625.79999999999
$ws.Range("M141").Value = -85445.79999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1862.75
$ws.Range("I32").Value = 1800.3954
$ws.Range("K32").Value = 1800.3954
$ws.Range("M32").Value = -1513.3954
$ws.Range("H45").Value = 2022.4762
$ws.Range("I45").Value = 2022.4762
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2022.4762
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1645.4762
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 41668364
$ws.Range("I61").Value = 45456172
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 45456172
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -45455960
$ws.Range("N61").Value = -2924
$ws.Range("H74").Value = 18203398
$ws.Range("I74").Value = 22248176
$ws.Range("K74").Value = 22248176
$ws.Range("M74").Value = -22247302
$ws.Range("H77").Value = 18203398
$ws.Range("I77").Value = 22248176
$ws.Range("K77").Value = 111240880
$ws.Range("M77").Value = -111236512
$ws.Range("H132").Value = 19638302
$ws.Range("I132").Value = 3254.439
$ws.Range("J132").Value = 100142000
$ws.Range("K132").Value = 9763.316999999999
$ws.Range("L132").Value = 300426000
$ws.Range("M132").Value = -7233.316999999999
$ws.Range("N132").Value = -300431060
$ws.Range("H136").Value = 41668364
$ws.Range("I136").Value = 45456172
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 136368516
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -136365966
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1331
$ws.Range("I20").Value = 926.4400000000001
$ws.Range("J20").Value = 2173.8333
$ws.Range("K20").Value = 926.4400000000001
$ws.Range("L20").Value = 2173.8333
$ws.Range("M20").Value = -679.4400000000001
$ws.Range("N20").Value = -2667.8333
$ws.Range("H26").Value = 27599.5
$ws.Range("I26").Value = 27599.5
$ws.Range("K26").Value = 27599.5
$ws.Range("M26").Value = -27307.5
$ws.Range("H28").Value = 60000
$ws.Range("J28").Value = 60000
$ws.Range("L28").Value = 60000
$ws.Range("N28").Value = -60588
$ws.Range("H64").Value = 1156
$ws.Range("I64").Value = 1200.1428
$ws.Range("J64").Value = 1078.75
$ws.Range("K64").Value = 1200.1428
$ws.Range("L64").Value = 1078.75
$ws.Range("M64").Value = -975.1428000000001
$ws.Range("N64").Value = -1528.75
$ws.Range("H67").Value = 1156
$ws.Range("I67").Value = 1200.1428
$ws.Range("J67").Value = 1078.75
$ws.Range("K67").Value = 1200.1428
$ws.Range("L67").Value = 1078.75
$ws.Range("M67").Value = -420.1428000000001
$ws.Range("N67").Value = -2638.75
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 2728
$ws.Range("I134").Value = 2538.6155
$ws.Range("J134").Value = 3343.5
$ws.Range("K134").Value = 7615.8465
$ws.Range("L134").Value = 10030.5
$ws.Range("M134").Value = -5080.8465
$ws.Range("N134").Value = -15100.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4365.6665
$ws.Range("I58").Value = 4821.2666
$ws.Range("J58").Value = 3606.3333
$ws.Range("K58").Value = 4821.2666
$ws.Range("L58").Value = 3606.3333
$ws.Range("M58").Value = -4618.2666
$ws.Range("N58").Value = -4012.3333
$ws.Range("H99").Value = 5451.6924
$ws.Range("I99").Value = 1200
$ws.Range("K99").Value = 1200
$ws.Range("M99").Value = 298
$ws.Range("H126").Value = 5451.6924
$ws.Range("I126").Value = 1200
$ws.Range("K126").Value = 3600
$ws.Range("M126").Value = -1130
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 4365.6665
$ws.Range("I136").Value = 4821.2666
$ws.Range("J136").Value = 3606.3333
$ws.Range("K136").Value = 14463.7998
$ws.Range("L136").Value = 10818.9999
$ws.Range("M136").Value = -11913.7998
$ws.Range("N136").Value = -15918.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1375287.8
$ws.Range("I32").Value = 200300
$ws.Range("J32").Value = 3333600.8
$ws.Range("K32").Value = 600900
$ws.Range("L32").Value = 10000802.4
$ws.Range("M32").Value = -600617
$ws.Range("N32").Value = -10001368.4
$ws.Range("H88").Value = 5996.25
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 5996.25
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H109").Value = 5303.727
$ws.Range("I109").Value = 6268.143
$ws.Range("J109").Value = 3616
$ws.Range("K109").Value = 18804.429
$ws.Range("L109").Value = 10848
$ws.Range("M109").Value = -17764.429
$ws.Range("N109").Value = -12928
$ws.Range("H131").Value = 45986.883
$ws.Range("J131").Value = 10323.214
$ws.Range("L131").Value = 30969.642
$ws.Range("N131").Value = -41049.642
$ws.Range("H133").Value = 6670.6665
$ws.Range("J133").Value = 19950
$ws.Range("L133").Value = 59850
$ws.Range("N133").Value = -69970
$ws.Range("H134").Value = 1678.3846
$ws.Range("J134").Value = 4891.2
$ws.Range("L134").Value = 14673.6
$ws.Range("N134").Value = -24813.6
$ws.Range("H138").Value = 3975
$ws.Range("J138").Value = 3033
$ws.Range("L138").Value = 9099
$ws.Range("N138").Value = -19379
$ws.Range("H139").Value = 2524.625
$ws.Range("I139").Value = 2029.9286
$ws.Range("J139").Value = 5987.5
$ws.Range("K139").Value = 6089.7858
$ws.Range("L139").Value = 17962.5
$ws.Range("M139").Value = -949.7857999999997
$ws.Range("N139").Value = -28242.5
$ws.Range("H141").Value = 7059.5713
$ws.Range("I141").Value = 3100.25
$ws.Range("J141").Value = 12338.667
$ws.Range("K141").Value = 9300.75
$ws.Range("L141").Value = 37016.001
$ws.Range("M141").Value = -4120.75
$ws.Range("N141").Value = -47376.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 11637.541
$ws.Range("I132").Value = 8029.4165
$ws.Range("J132").Value = 18298.691
$ws.Range("K132").Value = 24088.2495
$ws.Range("L132").Value = 54896.073
$ws.Range("M132").Value = -21558.2495
$ws.Range("N132").Value = -59956.073
$ws.Range("H137").Value = 63998
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H141").Value = 84999
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 84999
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 84999
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -95359

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3896.3333
$ws.Range("I61").Value = 2018.8
$ws.Range("K61").Value = 2018.8
$ws.Range("M61").Value = -1816.8
$ws.Range("H113").Value = 3896.3333
$ws.Range("I113").Value = 2018.8
$ws.Range("K113").Value = 2018.8
$ws.Range("M113").Value = 151.2
$ws.Range("H139").Value = 41905
$ws.Range("I139").Value = 31000
$ws.Range("K139").Value = 31000
$ws.Range("M139").Value = -25860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11342.714
$ws.Range("I62").Value = 14466.333
$ws.Range("K62").Value = 14466.333
$ws.Range("M62").Value = -13842.333
$ws.Range("H65").Value = 11342.714
$ws.Range("I65").Value = 14466.333
$ws.Range("K65").Value = 72331.66500000001
$ws.Range("M65").Value = -69211.66500000001
$ws.Range("H136").Value = 2071.3333
$ws.Range("J136").Value = 2014.1428
$ws.Range("L136").Value = 6042.428400000001
$ws.Range("N136").Value = -11142.4284
